$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C45").Value = 6218625100
$ws.Range("J45").Value = 6218625100
$ws.Range("K45").Value = 6218625101
$ws.Range("L45").Value = 3391977328
$ws.Range("M45").Value = 2826647773
$ws.Range("N45").Value = 47173352227
$ws.Range("C46").Value = 6218625100
$ws.Range("J46").Value = 6218625100
$ws.Range("K46").Value = 6218625101
$ws.Range("L46").Value = 1695988665
$ws.Range("M46").Value = 1695988663
$ws.Range("N46").Value = 28304011337
$ws.Range("C47").Value = 6218625100
$ws.Range("J47").Value = 6218625100
$ws.Range("K47").Value = 6218625101
$ws.Range("L47").Value = 2
$ws.Range("M47").Value = 1695988663
$ws.Range("N47").Value = 28304011337
$ws.Range("K48").Value = 2
$ws.Range("L48").Value = 2
$ws.Range("K49").Value = 2
$ws.Range("L49").Value = 2
$ws.Range("K50").Value = 2
$ws.Range("L50").Value = 2
$ws.Range("K51").Value = 2
$ws.Range("L51").Value = 2
$ws.Range("C52").Value = 8672784200
$ws.Range("J52").Value = 8672784200
$ws.Range("K52").Value = 8672784202
$ws.Range("L52").Value = 6349717006
$ws.Range("M52").Value = 2323067196
$ws.Range("N52").Value = 27676932804
$ws.Range("C53").Value = 8672784200
$ws.Range("J53").Value = 8672784200
$ws.Range("K53").Value = 8672784202
$ws.Range("L53").Value = 3639471943
$ws.Range("M53").Value = 2710245063
$ws.Range("N53").Value = 32289754937
$ws.Range("C54").Value = 8672784200
$ws.Range("J54").Value = 8672784200
$ws.Range("K54").Value = 8672784202
$ws.Range("L54").Value = 2
$ws.Range("M54").Value = 3639471941
$ws.Range("N54").Value = 43360528059
$ws.Range("K55").Value = 2
$ws.Range("L55").Value = 2
$ws.Range("K56").Value = 2
$ws.Range("L56").Value = 2
$ws.Range("K57").Value = 2
$ws.Range("L57").Value = 2
$ws.Range("C58").Value = 15637301200
$ws.Range("J58").Value = 15637301200
$ws.Range("K58").Value = 15637301202
$ws.Range("L58").Value = 15637301202
$ws.Range("C59").Value = 6097975400
$ws.Range("J59").Value = 6097975400
$ws.Range("K59").Value = 21735276602
$ws.Range("L59").Value = 21735276602
$ws.Range("K60").Value = 21735276602
$ws.Range("L60").Value = 21735276602
$ws.Range("C61").Value = 12700077000
$ws.Range("J61").Value = 12700077000
$ws.Range("K61").Value = 34435353602
$ws.Range("L61").Value = 34435353602
$ws.Range("K62").Value = 34435353602
$ws.Range("L62").Value = 34435353602
$ws.Range("K63").Value = 34435353602
$ws.Range("L63").Value = 34435353602
$ws.Range("K64").Value = 34435353602
$ws.Range("L64").Value = 34435353602
$ws.Range("K65").Value = 34435353602
$ws.Range("L65").Value = 34435353602
